$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal string into the cell without Excel's automatic
    # text->number coercion, and without leaving a residual number-format
    # style on the cell (matches original inlineStr cells with no "s" attr).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "36.486.03"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "1.937.80"
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "242.15"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$ws.Range("E6").Value = "  -1.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "56.42"
$ws.Range("E8").Value = "  -3.24%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.360"
$ws.Range("E9").Value = "  -3.33%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0808"
$ws.Range("E10").Value = "  -2.91%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.103"

# Row 12
$ws.Range("D12").Value = "2.221.87"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("E13").Value = "  -2.42%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.801"
$ws.Range("E14").Value = "  -3.79%  "

# Row 15
Set-TextValue $ws.Range("D15") "13.23"
$ws.Range("E15").Value = "  -2.53%  "

# Row 16
$ws.Range("E16").Value = "  -4.02%  "

# Row 17
$ws.Range("D17").Value = "1.934.83"
$ws.Range("E17").Value = "  -0.84%  "

# Row 18
$ws.Range("D18").Value = "36.421.69"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
Set-TextValue $ws.Range("D19") "68.96"
$ws.Range("E19").Value = "  -1.64%  "

# Row 20
$ws.Range("E20").Value = "  -2.91%  "

# Row 21
Set-TextValue $ws.Range("D21") "226.10"
$ws.Range("E21").Value = "  -1.89%  "

# Row 22
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("E23").Value = "  -0.21%  "

# Row 24
$ws.Range("E24").Value = "  -5.91%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.10"
$ws.Range("E26").Value = "  -4.26%  "

# Row 27
Set-TextValue $ws.Range("D27") "159.85"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28
$ws.Range("E28").Value = "  +8.51%  "

# Row 29
$ws.Range("E29").Value = "  -3.15%  "

# Row 30
$ws.Range("E30").Value = "  -0.91%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.08"
$ws.Range("E31").Value = "  -6.74%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.56"
$ws.Range("E32").Value = "  -3.61%  "

# Row 33
$ws.Range("E33").Value = "  -3.94%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.12"
$ws.Range("E34").Value = "  -4.48%  "

# Row 35
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("E36").Value = "  -0.82%  "

# Row 37
$ws.Range("E37").Value = "  -1.80%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.18"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39
$ws.Range("E39").Value = "  +9.82%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0984"
$ws.Range("E40").Value = "  -0.16%  "

# Row 41
Set-TextValue $ws.Range("D41") "2.91"
$ws.Range("E41").Value = "  +1.05%  "

# Row 42
$ws.Range("E42").Value = "  -1.41%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.14"
$ws.Range("E43").Value = "  -4.20%  "

# Row 44
Set-TextValue $ws.Range("D44") "15.61"
$ws.Range("E44").Value = "  -1.00%  "

# Row 45
$ws.Range("D45").Value = "1.332.29"
$ws.Range("E45").Value = "  -0.97%  "

# Row 46
$ws.Range("E46").Value = "  -1.94%  "

# Row 47
Set-TextValue $ws.Range("D47") "85.47"
$ws.Range("E47").Value = "  -3.79%  "

# Row 48
$ws.Range("E48").Value = "  -4.24%  "

# Row 49
$ws.Range("E49").Value = "  -0.61%  "

# Row 50
$ws.Range("D50").Value = "2.113.76"
$ws.Range("E50").Value = "  -0.93%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "42.84"
$ws.Range("E51").Value = "  -5.25%  "
